# update TestExcel & read data from excel
# Replace the single "HH" label in C3 with a small block of numeric test
# data (rows 3-5) that a "read data from excel" test would consume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: was just C3 = "HH" (shared string) -> now three numbers
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 3

# Row 4: new numeric data
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 5

# Row 5: new numeric data
$ws.Range("A5").Value = 6

# Move the active selection to G1, as in the saved workbook
[void]$ws.Range("G1").Select()
